$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 5 - pushes the existing rows 5-23
# (the product-list rows) down to become rows 6-24.
$ws.Rows.Item(5).Insert()

# Row insertion doesn't automatically re-point the conditional formatting
# ranges in this runtime, so do it explicitly so they keep covering the
# same logical rows (L1:Z<lastRow> and the product-code column A).
$cf1 = $ws.Range("L1:Z1048564").FormatConditions.Item(1)
$cf1.ModifyAppliesToRange($ws.Range("L1:Z1048565")) | Out-Null

$cf2 = $ws.Range("A9:A23").FormatConditions.Item(1)
$cf2.ModifyAppliesToRange($ws.Range("A10:A24")) | Out-Null

# Leave the selection on the newly inserted row.
$ws.Range("A5").Select() | Out-Null
